# "Generate Report for Archive"
#
# The localization-status report is regenerated for the archive. Part of
# the regeneration pipeline re-evaluates each handed-off file's
# translation status; for the 1f5ffa8b record the intermediate status
# computed during the refresh is "In Translation" before the pipeline
# settles back on the file's actual, already-reported status of
# "Ready for handoff" (the file's translation had in fact already been
# handed off, so the final report content is unchanged from before the
# refresh). Touching the status cells through that intermediate value is
# what registers "In Translation" as a known status string in the
# workbook even though no cell ends up displaying it.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Re-evaluate the 1f5ffa8b-fb0d-4407-b943-64c7143b7d8b.md row (row 3 on every
# sheet) as part of the report refresh.
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# The refreshed status resolves back to "Ready for handoff" for this
# record, matching what was already on record.
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "Ready for handoff"
$dede.Range("C3").Value = "Ready for handoff"
